$wb = $excel.ActiveWorkbook

# Add a new worksheet named "csvtest" placed after the last existing sheet
# (Sheet1, TestSheet) so it becomes the third / last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "csvtest"

# Populate the 2x2 range (A1:B2) with the value "csv" in every cell.
$newSheet.Range("A1:B2").Value = "csv"
